$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("company_list")

# Columns J (영업이익(비지배)) and O (자본총계(비지배)) are dropped entirely
# for rows 2-6, and AG/AH (현금DPS / 현금배당수익률) are dropped for row 6.
$ws.Range("J2:J6").ClearContents()
$ws.Range("O2:O6").ClearContents()
$ws.Range("AG6:AH6").ClearContents()

# Rows 7-9 lose all their financial data (columns D through AJ); only
# the leading A/B/C (row#, period, year label) columns remain.
$ws.Range("D7:AJ9").ClearContents()

# Row 2 (2014/12 IFRS연결)
$ws.Range("D2").Value = 2542
$ws.Range("E2").Value = -134
$ws.Range("F2").Value = -144
$ws.Range("G2").Value = -1623
$ws.Range("H2").Value = -1245
$ws.Range("I2").Value = -1245
$ws.Range("K2").Value = 3659
$ws.Range("L2").Value = 2556
$ws.Range("M2").Value = 1102
$ws.Range("N2").Value = 1102
$ws.Range("P2").Value = 911
$ws.Range("Q2").Value = 275
$ws.Range("R2").Value = 1019
$ws.Range("S2").Value = -1246
$ws.Range("T2").Value = 26
$ws.Range("U2").Value = 248
$ws.Range("V2").Value = 1294
$ws.Range("W2").Value = -5.28
$ws.Range("X2").Value = -48.98
$ws.Range("Y2").Value = -69.86
$ws.Range("Z2").Value = -24.8
$ws.Range("AA2").Value = 231.84
$ws.Range("AB2").Value = 47.92
$ws.Range("AC2").Value = -684
$ws.Range("AD2").Value = -0.37
$ws.Range("AE2").Value = 641
$ws.Range("AF2").Value = 0.39
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 182126320

# Row 3 (2015/12 IFRS연결)
$ws.Range("D3").Value = 2037
$ws.Range("E3").Value = 98
$ws.Range("F3").Value = 98
$ws.Range("G3").Value = -342
$ws.Range("H3").Value = 132
$ws.Range("I3").Value = 132
$ws.Range("K3").Value = 2762
$ws.Range("L3").Value = 1104
$ws.Range("M3").Value = 1658
$ws.Range("N3").Value = 1658
$ws.Range("P3").Value = 911
$ws.Range("Q3").Value = -131
$ws.Range("R3").Value = 1046
$ws.Range("S3").Value = -864
$ws.Range("T3").Value = 26
$ws.Range("U3").Value = -157
$ws.Range("V3").Value = 425
$ws.Range("W3").Value = 4.79
$ws.Range("X3").Value = 6.48
$ws.Range("Y3").Value = 9.550000000000001
$ws.Range("Z3").Value = 4.11
$ws.Range("AA3").Value = 66.56
$ws.Range("AB3").Value = 62.39
$ws.Range("AC3").Value = 72
$ws.Range("AD3").Value = 6.84
$ws.Range("AE3").Value = 964
$ws.Range("AF3").Value = 0.51
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 182126320

# Row 4 (2016/12 IFRS연결)
$ws.Range("D4").Value = 2048
$ws.Range("E4").Value = 109
$ws.Range("F4").Value = 109
$ws.Range("G4").Value = -143
$ws.Range("H4").Value = -207
$ws.Range("I4").Value = -207
$ws.Range("K4").Value = 2150
$ws.Range("L4").Value = 634
$ws.Range("M4").Value = 1515
$ws.Range("N4").Value = 1515
$ws.Range("P4").Value = 911
$ws.Range("Q4").Value = 144
$ws.Range("R4").Value = 69
$ws.Range("S4").Value = -186
$ws.Range("T4").Value = 17
$ws.Range("U4").Value = 126
$ws.Range("V4").Value = 241
$ws.Range("W4").Value = 5.32
$ws.Range("X4").Value = -10.11
$ws.Range("Y4").Value = -13.05
$ws.Range("Z4").Value = -8.43
$ws.Range("AA4").Value = 41.85
$ws.Range("AB4").Value = 39.1
$ws.Range("AC4").Value = -114
$ws.Range("AD4").Value = -6.45
$ws.Range("AE4").Value = 881
$ws.Range("AF4").Value = 0.83
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 182126320

# Row 5 (2017/12 IFRS연결)
$ws.Range("D5").Value = 1968
$ws.Range("E5").Value = 38
$ws.Range("F5").Value = 38
$ws.Range("G5").Value = -130
$ws.Range("H5").Value = -10
$ws.Range("I5").Value = -10
$ws.Range("K5").Value = 1975
$ws.Range("L5").Value = 584
$ws.Range("M5").Value = 1391
$ws.Range("N5").Value = 1391
$ws.Range("P5").Value = 949
$ws.Range("Q5").Value = 80
$ws.Range("R5").Value = 18
$ws.Range("S5").Value = -83
$ws.Range("T5").Value = 8
$ws.Range("U5").Value = 72
$ws.Range("V5").Value = 143
$ws.Range("W5").Value = 1.96
$ws.Range("X5").Value = -0.49
$ws.Range("Y5").Value = -0.66
$ws.Range("Z5").Value = -0.46
$ws.Range("AA5").Value = 41.95
$ws.Range("AB5").Value = 37.91
$ws.Range("AC5").Value = -5
$ws.Range("AD5").Value = -125.9
$ws.Range("AE5").Value = 775
$ws.Range("AF5").Value = 0.85
$ws.Range("AG5").Value = 0
$ws.Range("AH5").Value = 0
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 189745366

# Row 6 (2018/12 IFRS연결)
$ws.Range("D6").Value = 2371
$ws.Range("E6").Value = 16
$ws.Range("F6").Value = 16
$ws.Range("G6").Value = -47
$ws.Range("H6").Value = -40
$ws.Range("I6").Value = -40
$ws.Range("K6").Value = 2884
$ws.Range("L6").Value = 1508
$ws.Range("M6").Value = 1376
$ws.Range("N6").Value = 1376
$ws.Range("P6").Value = 1006
$ws.Range("Q6").Value = 68
$ws.Range("R6").Value = -848
$ws.Range("S6").Value = 928
$ws.Range("T6").Value = 4
$ws.Range("U6").Value = 64
$ws.Range("V6").Value = 1016
$ws.Range("W6").Value = 0.66
$ws.Range("X6").Value = -1.68
$ws.Range("Y6").Value = -2.87
$ws.Range("Z6").Value = -1.64
$ws.Range("AA6").Value = 109.65
$ws.Range("AB6").Value = 184.29
$ws.Range("AC6").Value = -20
$ws.Range("AD6").Value = -41.55
$ws.Range("AE6").Value = 720
$ws.Range("AF6").Value = 1.15
$ws.Range("AI6").Value = 0
$ws.Range("AJ6").Value = 201173933
